$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells (Price column D) keep their original text
# representation instead of being auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.402.35"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.45"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.70"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E7").Value = "  +1.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.78"
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07649"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.150"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.007"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.939"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.571.44"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.20"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.75"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.230"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.02"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.393.72"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.396"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.653"
$ws.Range("E26").Value = "  -10.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.18"
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "147.07"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.029"
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.74"
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.746.08"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.163"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.011"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9818"
$ws.Range("E34").Value = "  -5.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.12"
$ws.Range("E35").Value = "  -1.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08477"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02537"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.378"
$ws.Range("E38").Value = "  +11.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2316"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06552"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.422"
$ws.Range("E41").Value = "  -2.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.46"
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6386"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.06"
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.805"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5975"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.290"
$ws.Range("E48").Value = "  +2.93%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.095"
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.75"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07330"
